$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 774587.9
$ws.Range("I9").Value = 865665.9
$ws.Range("J9").Value = 425
$ws.Range("K9").Value = 865665.9
$ws.Range("L9").Value = 425
$ws.Range("M9").Value = -865496.9
$ws.Range("N9").Value = -763
$ws.Range("H11").Value = 311.55554
$ws.Range("I11").Value = 311.55554
$ws.Range("K11").Value = 311.55554
$ws.Range("M11").Value = -171.55554
$ws.Range("H18").Value = 125001304
$ws.Range("I18").Value = 927.1667
$ws.Range("K18").Value = 927.1667
$ws.Range("M18").Value = -643.1667
$ws.Range("H29").Value = 1595.6666
$ws.Range("J29").Value = 1895
$ws.Range("L29").Value = 5685
$ws.Range("N29").Value = -6247
$ws.Range("H33").Value = 394.5
$ws.Range("H51").Value = 83339000
$ws.Range("J51").Value = 5999.5
$ws.Range("L51").Value = 5999.5
$ws.Range("N51").Value = -6967.5
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("H64").Value = 83338620
$ws.Range("I64").Value = 5985
$ws.Range("J64").Value = 142861940
$ws.Range("K64").Value = 5985
$ws.Range("L64").Value = 142861940
$ws.Range("M64").Value = -5737
$ws.Range("N64").Value = -142862436
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
$ws.Range("H67").Value = 83338620
$ws.Range("I67").Value = 5985
$ws.Range("J67").Value = 142861940
$ws.Range("K67").Value = 5985
$ws.Range("L67").Value = 142861940
$ws.Range("M67").Value = -5127
$ws.Range("N67").Value = -142863656
$ws.Range("H74").Value = 12142.423
$ws.Range("I74").Value = 14464.294
$ws.Range("K74").Value = 14464.294
$ws.Range("M74").Value = -13528.294
$ws.Range("H77").Value = 12142.423
$ws.Range("I77").Value = 14464.294
$ws.Range("K77").Value = 72321.47
$ws.Range("M77").Value = -67641.47
$ws.Range("H88").Value = 4006646.5
$ws.Range("I88").Value = 20003058
$ws.Range("J88").Value = 7543.7
$ws.Range("K88").Value = 20003058
$ws.Range("L88").Value = 7543.7
$ws.Range("M88").Value = -20002652
$ws.Range("N88").Value = -8355.700000000001
$ws.Range("H91").Value = 4006646.5
$ws.Range("I91").Value = 20003058
$ws.Range("J91").Value = 7543.7
$ws.Range("K91").Value = 20003058
$ws.Range("L91").Value = 7543.7
$ws.Range("M91").Value = -20001654
$ws.Range("N91").Value = -10351.7
$ws.Range("H92").Value = 509.5
$ws.Range("I92").Value = 527.8333
$ws.Range("J92").Value = 399.5
$ws.Range("K92").Value = 527.8333
$ws.Range("L92").Value = 399.5
$ws.Range("M92").Value = 720.1667
$ws.Range("N92").Value = -2895.5
$ws.Range("H98").Value = 2663.9736
$ws.Range("I98").Value = 2632.8333
$ws.Range("K98").Value = 2632.8333
$ws.Range("M98").Value = -1134.8333
$ws.Range("H100").Value = 5899.2
$ws.Range("I100").Value = 5899.2
$ws.Range("K100").Value = 5899.2
$ws.Range("M100").Value = -5358.2
$ws.Range("H112").Value = 2771.6667
$ws.Range("I112").Value = 666.6667
$ws.Range("K112").Value = 2000.0001
$ws.Range("M112").Value = -892.0001
$ws.Range("H118").Value = 599.8333
$ws.Range("I118").Value = 599.8333
$ws.Range("K118").Value = 1799.4999
$ws.Range("M118").Value = -142.4999
$ws.Range("H122").Value = 2663.9736
$ws.Range("I122").Value = 2632.8333
$ws.Range("K122").Value = 7898.499899999999
$ws.Range("M122").Value = -5448.499899999999
$ws.Range("H125").Value = 601.5714
$ws.Range("J125").Value = 524.75
$ws.Range("L125").Value = 4722.75
$ws.Range("N125").Value = -9642.75
$ws.Range("H132").Value = 2619664.2
$ws.Range("I132").Value = 5592
$ws.Range("K132").Value = 16776
$ws.Range("M132").Value = -14246
$ws.Range("H135").Value = 623.2973
$ws.Range("J135").Value = 1245.1666
$ws.Range("L135").Value = 11206.4994
$ws.Range("N135").Value = -16276.4994
$ws.Range("H137").Value = 1138818.1
$ws.Range("I137").Value = 1353387.1
$ws.Range("J137").Value = 4667.7144
$ws.Range("K137").Value = 4060161.3
$ws.Range("L137").Value = 14003.1432
$ws.Range("M137").Value = -4057611.3
$ws.Range("N137").Value = -19103.1432
$ws.Range("H138").Value = 2705.9622
$ws.Range("I138").Value = 1299.25
$ws.Range("J138").Value = 3870.138
$ws.Range("K138").Value = 3897.75
$ws.Range("L138").Value = 11610.414
$ws.Range("M138").Value = 1242.25
$ws.Range("N138").Value = -21890.414
$ws.Range("H141").Value = 2788
$ws.Range("I141").Value = 2557
$ws.Range("J141").Value = 3250
$ws.Range("K141").Value = 7671
$ws.Range("L141").Value = 9750
$ws.Range("M141").Value = -2491
$ws.Range("N141").Value = -20110

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 699
$ws.Range("I4").Value = 699
$ws.Range("K4").Value = 699
$ws.Range("M4").Value = -583
$ws.Range("H5").Value = 245.33333
$ws.Range("I5").Value = 245.33333
$ws.Range("K5").Value = 245.33333
$ws.Range("M5").Value = -133.33333
$ws.Range("H10").Value = 20004
$ws.Range("I10").Value = 20004
$ws.Range("K10").Value = 20004
$ws.Range("M10").Value = -19834
$ws.Range("H32").Value = 1674.0714
$ws.Range("I32").Value = 1648.9615
$ws.Range("J32").Value = 2000.5
$ws.Range("K32").Value = 1648.9615
$ws.Range("L32").Value = 2000.5
$ws.Range("M32").Value = -1361.9615
$ws.Range("N32").Value = -2574.5
$ws.Range("H61").Value = 2016.3684
$ws.Range("I61").Value = 1367.0333
$ws.Range("K61").Value = 1367.0333
$ws.Range("M61").Value = -1155.0333
$ws.Range("H74").Value = 110707.18
$ws.Range("I74").Value = 130606.19
$ws.Range("K74").Value = 130606.19
$ws.Range("M74").Value = -129732.19
$ws.Range("H77").Value = 110707.18
$ws.Range("I77").Value = 130606.19
$ws.Range("K77").Value = 653030.95
$ws.Range("M77").Value = -648662.95
$ws.Range("H102").Value = 4535.154
$ws.Range("I102").Value = 3945
$ws.Range("K102").Value = 3945
$ws.Range("M102").Value = -2323
$ws.Range("H110").Value = 737.5
$ws.Range("I110").Value = 737.5
$ws.Range("K110").Value = 737.5
$ws.Range("M110").Value = 1307.5
$ws.Range("H132").Value = 6669152
$ws.Range("I132").Value = 1972.8918
$ws.Range("K132").Value = 5918.6754
$ws.Range("M132").Value = -3388.6754
$ws.Range("H136").Value = 2016.3684
$ws.Range("I136").Value = 1367.0333
$ws.Range("K136").Value = 4101.0999
$ws.Range("M136").Value = -1551.0999
$ws.Range("H139").Value = 80730
$ws.Range("J139").Value = 80730
$ws.Range("L139").Value = 80730
$ws.Range("N139").Value = -91010

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 245.33333
$ws.Range("I4").Value = 245.33333
$ws.Range("K4").Value = 245.33333
$ws.Range("M4").Value = -130.33333
$ws.Range("H39").Value = 10000
$ws.Range("J39").Value = 10000
$ws.Range("L39").Value = 10000
$ws.Range("N39").Value = -10778
$ws.Range("H55").Value = 22821
$ws.Range("J55").Value = 26926.334
$ws.Range("L55").Value = 26926.334
$ws.Range("N55").Value = -27472.334
$ws.Range("H80").Value = 649.1429000000001
$ws.Range("I80").Value = 594.5
$ws.Range("J80").Value = 671
$ws.Range("K80").Value = 594.5
$ws.Range("L80").Value = 671
$ws.Range("M80").Value = 403.5
$ws.Range("N80").Value = -2667
$ws.Range("H82").Value = 41074.9
$ws.Range("I82").Value = 25757.857
$ws.Range("J82").Value = 76814.664
$ws.Range("K82").Value = 25757.857
$ws.Range("L82").Value = 76814.664
$ws.Range("M82").Value = -25374.857
$ws.Range("N82").Value = -77580.664
$ws.Range("H83").Value = 649.1429000000001
$ws.Range("I83").Value = 594.5
$ws.Range("J83").Value = 671
$ws.Range("K83").Value = 2972.5
$ws.Range("L83").Value = 3355
$ws.Range("M83").Value = 2019.5
$ws.Range("N83").Value = -13339
$ws.Range("H85").Value = 41074.9
$ws.Range("I85").Value = 25757.857
$ws.Range("J85").Value = 76814.664
$ws.Range("K85").Value = 25757.857
$ws.Range("L85").Value = 76814.664
$ws.Range("M85").Value = -24431.857
$ws.Range("N85").Value = -79466.664
$ws.Range("H86").Value = 4012.76
$ws.Range("I86").Value = 3690.7646
$ws.Range("J86").Value = 4697
$ws.Range("K86").Value = 3690.7646
$ws.Range("L86").Value = 4697
$ws.Range("M86").Value = -2567.7646
$ws.Range("N86").Value = -6943
$ws.Range("H89").Value = 4012.76
$ws.Range("I89").Value = 3690.7646
$ws.Range("J89").Value = 4697
$ws.Range("K89").Value = 18453.823
$ws.Range("L89").Value = 23485
$ws.Range("M89").Value = -12837.823
$ws.Range("N89").Value = -34717
$ws.Range("H94").Value = 57144670
$ws.Range("I94").Value = 71430170
$ws.Range("K94").Value = 71430170
$ws.Range("M94").Value = -71429719
$ws.Range("H105").Value = 10002865
$ws.Range("I105").Value = 529110.8
$ws.Range("J105").Value = 35717340
$ws.Range("K105").Value = 529110.8
$ws.Range("L105").Value = 35717340
$ws.Range("M105").Value = -527363.8
$ws.Range("N105").Value = -35720834
$ws.Range("H116").Value = 79742
$ws.Range("J116").Value = 79742
$ws.Range("L116").Value = 79742
$ws.Range("N116").Value = -88920
$ws.Range("H134").Value = 3579.7
$ws.Range("I134").Value = 2988.7932
$ws.Range("J134").Value = 5137.5454
$ws.Range("K134").Value = 8966.3796
$ws.Range("L134").Value = 15412.6362
$ws.Range("M134").Value = -6431.3796
$ws.Range("N134").Value = -20482.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 60339.5
$ws.Range("J18").Value = 60339.5
$ws.Range("L18").Value = 60339.5
$ws.Range("N18").Value = -60799.5
$ws.Range("H31").Value = 3681044.8
$ws.Range("I31").Value = 3369.9614
$ws.Range("J31").Value = 15633487
$ws.Range("K31").Value = 3369.9614
$ws.Range("L31").Value = 15633487
$ws.Range("M31").Value = -3074.9614
$ws.Range("N31").Value = -15634077
$ws.Range("H34").Value = 3681044.8
$ws.Range("I34").Value = 3369.9614
$ws.Range("J34").Value = 15633487
$ws.Range("K34").Value = 3369.9614
$ws.Range("L34").Value = 15633487
$ws.Range("M34").Value = -3167.9614
$ws.Range("N34").Value = -15633891
$ws.Range("H43").Value = 18968
$ws.Range("J43").Value = 20653.5
$ws.Range("L43").Value = 20653.5
$ws.Range("N43").Value = -21021.5
$ws.Range("H58").Value = 2446.7144
$ws.Range("I58").Value = 2165.8823
$ws.Range("K58").Value = 2165.8823
$ws.Range("M58").Value = -1962.8823
$ws.Range("H62").Value = 10004001
$ws.Range("I62").Value = 12504001
$ws.Range("K62").Value = 12504001
$ws.Range("M62").Value = -12503377
$ws.Range("H65").Value = 10004001
$ws.Range("I65").Value = 12504001
$ws.Range("K65").Value = 62520005
$ws.Range("M65").Value = -62516885
$ws.Range("H101").Value = 18968
$ws.Range("J101").Value = 20653.5
$ws.Range("L101").Value = 20653.5
$ws.Range("N101").Value = -27143.5
$ws.Range("H106").Value = 40203.332
$ws.Range("J106").Value = 45000
$ws.Range("L106").Value = 45000
$ws.Range("N106").Value = -47524
$ws.Range("H107").Value = 1724761
$ws.Range("I107").Value = 2632021.8
$ws.Range("J107").Value = 965.8
$ws.Range("K107").Value = 2632021.8
$ws.Range("L107").Value = 965.8
$ws.Range("M107").Value = -2630101.8
$ws.Range("N107").Value = -4805.8
$ws.Range("H122").Value = 2973.2104
$ws.Range("I122").Value = 2275.6365
$ws.Range("K122").Value = 6826.9095
$ws.Range("M122").Value = -4376.9095
$ws.Range("H124").Value = 49662.75
$ws.Range("J124").Value = 49662.75
$ws.Range("L124").Value = 49662.75
$ws.Range("N124").Value = -54572.75
$ws.Range("H132").Value = 7577376
$ws.Range("I132").Value = 1332.2727
$ws.Range("K132").Value = 3996.8181
$ws.Range("M132").Value = -1466.8181
$ws.Range("H134").Value = 3440.2258
$ws.Range("I134").Value = 2610.8262
$ws.Range("K134").Value = 7832.4786
$ws.Range("M134").Value = -5297.4786
$ws.Range("H136").Value = 2446.7144
$ws.Range("I136").Value = 2165.8823
$ws.Range("K136").Value = 6497.646900000001
$ws.Range("M136").Value = -3947.646900000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 434.63635
$ws.Range("I14").Value = 434.63635
$ws.Range("K14").Value = 1303.90905
$ws.Range("M14").Value = -1130.90905
$ws.Range("H15").Value = 604
$ws.Range("I15").Value = 85.333336
$ws.Range("J15").Value = 1771
$ws.Range("K15").Value = 256.000008
$ws.Range("L15").Value = 5313
$ws.Range("M15").Value = -116.000008
$ws.Range("N15").Value = -5593
$ws.Range("H34").Value = 1656.2778
$ws.Range("J34").Value = 1894.4
$ws.Range("L34").Value = 5683.200000000001
$ws.Range("N34").Value = -5851.200000000001
$ws.Range("H39").Value = 4240.091
$ws.Range("I39").Value = 3932.3333
$ws.Range("J39").Value = 4355.5
$ws.Range("K39").Value = 11796.9999
$ws.Range("L39").Value = 13066.5
$ws.Range("M39").Value = -11502.9999
$ws.Range("N39").Value = -13654.5
$ws.Range("H55").Value = 6313.1665
$ws.Range("J55").Value = 7498.3
$ws.Range("L55").Value = 22494.9
$ws.Range("N55").Value = -22848.9
$ws.Range("H86").Value = 1484.2
$ws.Range("I86").Value = 624
$ws.Range("J86").Value = 2057.6667
$ws.Range("K86").Value = 1872
$ws.Range("L86").Value = 6173.000100000001
$ws.Range("M86").Value = -686
$ws.Range("N86").Value = -8545.000100000001
$ws.Range("H89").Value = 1484.2
$ws.Range("I89").Value = 624
$ws.Range("J89").Value = 2057.6667
$ws.Range("K89").Value = 5616
$ws.Range("L89").Value = 18519.0003
$ws.Range("M89").Value = 312
$ws.Range("N89").Value = -30375.0003
$ws.Range("H93").Value = 5949
$ws.Range("J93").Value = 5949
$ws.Range("L93").Value = 17847
$ws.Range("N93").Value = -21591
$ws.Range("H100").Value = 2000
$ws.Range("J100").Value = 2000
$ws.Range("L100").Value = 6000
$ws.Range("N100").Value = -7622
$ws.Range("H116").Value = 4799.6665
$ws.Range("J116").Value = 5999.6665
$ws.Range("L116").Value = 17998.9995
$ws.Range("N116").Value = -24882.9995
$ws.Range("H118").Value = 16785.3
$ws.Range("I118").Value = 18481.75
$ws.Range("K118").Value = 55445.25
$ws.Range("M118").Value = -54202.25
$ws.Range("H122").Value = 1434.7941
$ws.Range("J122").Value = 1756.2693
$ws.Range("L122").Value = 15806.4237
$ws.Range("N122").Value = -20706.4237
$ws.Range("H132").Value = 26915.3
$ws.Range("I132").Value = 29617
$ws.Range("J132").Value = 2600
$ws.Range("K132").Value = 266553
$ws.Range("L132").Value = 23400
$ws.Range("M132").Value = -264023
$ws.Range("N132").Value = -28460
$ws.Range("H139").Value = 6406.276
$ws.Range("I139").Value = 14459.223
$ws.Range("J139").Value = 2782.45
$ws.Range("K139").Value = 43377.669
$ws.Range("L139").Value = 8347.349999999999
$ws.Range("M139").Value = -38237.669
$ws.Range("N139").Value = -18627.35
$ws.Range("H140").Value = 9261.314
$ws.Range("I140").Value = 6631.625
$ws.Range("J140").Value = 14998.818
$ws.Range("K140").Value = 19894.875
$ws.Range("L140").Value = 44996.454
$ws.Range("M140").Value = -14714.875
$ws.Range("N140").Value = -55356.454
$ws.Range("H141").Value = 4706.6875
$ws.Range("I141").Value = 4706.6875
$ws.Range("K141").Value = 14120.0625
$ws.Range("M141").Value = -8940.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5894503.5
$ws.Range("I11").Value = 12146257
$ws.Range("J11").Value = 59533.4
$ws.Range("K11").Value = 12146257
$ws.Range("L11").Value = 59533.4
$ws.Range("M11").Value = -12146118
$ws.Range("N11").Value = -59811.4
$ws.Range("H14").Value = 5148498
$ws.Range("I14").Value = 6004498
$ws.Range("K14").Value = 6004498
$ws.Range("M14").Value = -6004330
$ws.Range("H58").Value = 23999
$ws.Range("I58").Value = 22665.666
$ws.Range("K58").Value = 22665.666
$ws.Range("M58").Value = -22388.666
$ws.Range("H80").Value = 58826650
$ws.Range("I80").Value = 100002190
$ws.Range("K80").Value = 100002190
$ws.Range("M80").Value = -100001192
$ws.Range("H83").Value = 58826650
$ws.Range("I83").Value = 100002190
$ws.Range("K83").Value = 500010950
$ws.Range("M83").Value = -500005958
$ws.Range("H102").Value = 2004.0667
$ws.Range("I102").Value = 1734.4615
$ws.Range("K102").Value = 1734.4615
$ws.Range("M102").Value = -112.4614999999999
$ws.Range("H122").Value = 3729.0667
$ws.Range("I122").Value = 1603.7
$ws.Range("J122").Value = 7979.8
$ws.Range("K122").Value = 4811.1
$ws.Range("L122").Value = 23939.4
$ws.Range("M122").Value = -2361.1
$ws.Range("N122").Value = -28839.4
$ws.Range("H126").Value = 88900410
$ws.Range("I126").Value = 7495.5
$ws.Range("J126").Value = 148162350
$ws.Range("K126").Value = 22486.5
$ws.Range("L126").Value = 444487050
$ws.Range("M126").Value = -20016.5
$ws.Range("N126").Value = -444491990
$ws.Range("H132").Value = 2476.7827
$ws.Range("I132").Value = 2093.6875
$ws.Range("J132").Value = 3352.4285
$ws.Range("K132").Value = 6281.0625
$ws.Range("L132").Value = 10057.2855
$ws.Range("M132").Value = -3751.0625
$ws.Range("N132").Value = -15117.2855
$ws.Range("H136").Value = 35666.5
$ws.Range("J136").Value = 35666.5
$ws.Range("L136").Value = 106999.5
$ws.Range("N136").Value = -112099.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3514.55
$ws.Range("I40").Value = 3739.3333
$ws.Range("K40").Value = 3739.3333
$ws.Range("M40").Value = -3603.3333
$ws.Range("H46").Value = 2000
$ws.Range("I46").Value = 2000
$ws.Range("K46").Value = 2000
$ws.Range("M46").Value = -1812
$ws.Range("H61").Value = 8306.9375
$ws.Range("J61").Value = 37334
$ws.Range("L61").Value = 37334
$ws.Range("N61").Value = -37738
$ws.Range("H64").Value = 33356.168
$ws.Range("J64").Value = 33356.168
$ws.Range("L64").Value = 33356.168
$ws.Range("N64").Value = -33806.168
$ws.Range("H67").Value = 33356.168
$ws.Range("J67").Value = 33356.168
$ws.Range("L67").Value = 33356.168
$ws.Range("N67").Value = -34916.168
$ws.Range("H68").Value = 2600
$ws.Range("J68").Value = 2692.4614
$ws.Range("L68").Value = 2692.4614
$ws.Range("N68").Value = -4190.4614
$ws.Range("H71").Value = 2600
$ws.Range("J71").Value = 2692.4614
$ws.Range("L71").Value = 13462.307
$ws.Range("N71").Value = -20950.307
$ws.Range("H82").Value = 1521
$ws.Range("I82").Value = 1476.3636
$ws.Range("J82").Value = 2503
$ws.Range("K82").Value = 1476.3636
$ws.Range("L82").Value = 2503
$ws.Range("M82").Value = -1115.3636
$ws.Range("N82").Value = -3225
$ws.Range("H85").Value = 1521
$ws.Range("I85").Value = 1476.3636
$ws.Range("J85").Value = 2503
$ws.Range("K85").Value = 1476.3636
$ws.Range("L85").Value = 2503
$ws.Range("M85").Value = -228.3635999999999
$ws.Range("N85").Value = -4999
$ws.Range("H100").Value = 3997.1428
$ws.Range("I100").Value = 3996.8333
$ws.Range("K100").Value = 3996.8333
$ws.Range("M100").Value = -3455.8333
$ws.Range("H113").Value = 8306.9375
$ws.Range("J113").Value = 37334
$ws.Range("L113").Value = 37334
$ws.Range("N113").Value = -41674
$ws.Range("H132").Value = 3835.1177
$ws.Range("I132").Value = 3063.818
$ws.Range("J132").Value = 5249.1665
$ws.Range("K132").Value = 9191.454000000002
$ws.Range("L132").Value = 15747.4995
$ws.Range("M132").Value = -6661.454000000002
$ws.Range("N132").Value = -20807.4995
$ws.Range("H136").Value = 3943.1875
$ws.Range("I136").Value = 2720.7307
$ws.Range("J136").Value = 9240.5
$ws.Range("K136").Value = 8162.1921
$ws.Range("L136").Value = 27721.5
$ws.Range("M136").Value = -5612.1921
$ws.Range("N136").Value = -32821.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 13333.333
$ws.Range("I18").Value = 10000
$ws.Range("K18").Value = 10000
$ws.Range("M18").Value = -9827
$ws.Range("H58").Value = 14518.75
$ws.Range("I58").Value = 14518.75
$ws.Range("K58").Value = 14518.75
$ws.Range("M58").Value = -14210.75
$ws.Range("H63").Value = 26749.666
$ws.Range("J63").Value = 26749.666
$ws.Range("L63").Value = 26749.666
$ws.Range("N63").Value = -27997.666
$ws.Range("H66").Value = 26749.666
$ws.Range("J66").Value = 26749.666
$ws.Range("L66").Value = 80248.99800000001
$ws.Range("N66").Value = -86488.99800000001
$ws.Range("H81").Value = 6654.5835
$ws.Range("J81").Value = 5419.375
$ws.Range("L81").Value = 10838.75
$ws.Range("N81").Value = -12960.75
$ws.Range("H84").Value = 6654.5835
$ws.Range("J84").Value = 5419.375
$ws.Range("L84").Value = 54193.75
$ws.Range("N84").Value = -64801.75
$ws.Range("H100").Value = 55556170
$ws.Range("I100").Value = 573.3077
$ws.Range("J100").Value = 200000720
$ws.Range("K100").Value = 1146.6154
$ws.Range("L100").Value = 400001440
$ws.Range("M100").Value = -605.6153999999999
$ws.Range("N100").Value = -400002522
$ws.Range("H107").Value = 1066.4
$ws.Range("I107").Value = 919.6667
$ws.Range("J107").Value = 1286.5
$ws.Range("K107").Value = 2759.0001
$ws.Range("L107").Value = 3859.5
$ws.Range("M107").Value = -839.0001000000002
$ws.Range("N107").Value = -7699.5
$ws.Range("H113").Value = 790.7917
$ws.Range("I113").Value = 684.8125
$ws.Range("K113").Value = 2054.4375
$ws.Range("M113").Value = 115.5625
$ws.Range("H122").Value = 17859728
$ws.Range("I122").Value = 2979.182
$ws.Range("J122").Value = 83334470
$ws.Range("K122").Value = 8937.545999999998
$ws.Range("L122").Value = 250003410
$ws.Range("M122").Value = -6487.545999999998
$ws.Range("N122").Value = -250008310
$ws.Range("H132").Value = 2052.2942
$ws.Range("I132").Value = 2099.372
$ws.Range("J132").Value = 1799.25
$ws.Range("K132").Value = 6298.116
$ws.Range("L132").Value = 5397.75
$ws.Range("M132").Value = -3768.116
$ws.Range("N132").Value = -10457.75
$ws.Range("H133").Value = 86141.86
$ws.Range("J133").Value = 86141.86
$ws.Range("L133").Value = 86141.86
$ws.Range("N133").Value = -96261.86
$ws.Range("H136").Value = 8349.025
$ws.Range("I136").Value = 8800.028
$ws.Range("K136").Value = 26400.084
$ws.Range("M136").Value = -23850.084
$ws.Range("H140").Value = 74196
$ws.Range("J140").Value = 74196
$ws.Range("L140").Value = 74196
$ws.Range("N140").Value = -84556
$ws.Range("H141").Value = 69999
$ws.Range("J141").Value = 69999
$ws.Range("L141").Value = 69999
$ws.Range("N141").Value = -80359

Write-Output "Applied 597 cell updates across 8 sheets"